# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.873.83'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '3.389.58'
$ws.Range("E3").Value = '  -2.90%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.45%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '3.388.58'
$ws.Range("E9").Value = '  -2.80%  '
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.93'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.412'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '3.976.85'
$ws.Range("E13").Value = '  -2.55%  '
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("D16").Value = '65.993.83'
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = '3.393.75'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.527'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.83%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.47%  '
$ws.Range("E35").Value = '  -3.53%  '
$ws.Range("E36").Value = '  -2.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.858'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.60'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.664.73'
$ws.Range("E43").Value = '  -4.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0677'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.90%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '333.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.51%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '24.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0284'
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = '  +2.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.60%  '
